# Refresh the cryptos list (Sheet1) with the latest scraped Price (D) and
# Volume(1h) (E) text for rows 2-51, mirroring the source GitHub Actions job.
#
# Price/Volume are stored as literal text (e.g. "1.000", "27.170.62",
# "  +0.80%  ") rather than numbers, so plain `.Value = "<numeric-looking>"`
# assignments would let Excel silently coerce them into numbers (dropping
# trailing zeros, merging "27.170.62" into 27170.62, etc.). For the Price
# column we therefore force a Text format before writing the value and then
# clear the formatting again afterwards so the cell keeps its original
# (unstyled) look while the text itself stays exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.170.62'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.908.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5248'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.11%  '
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07279'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +2.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8989'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07677'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.886.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.75'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.252'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008554'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("E18").Value = '  +2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.234.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.089'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.128.79'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.447'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.318'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +10.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.86'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.730'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.83'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.968'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.02%  '
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09215'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05077'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.246'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7796'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.993'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.310'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.610'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5681'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01996'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.073'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.030'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.633'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.70'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.81%  '
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.20'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.83%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.606'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.49'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.23'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.58%  '
